# Applies the "ENTREGA 5 / ENTREGA 6" rows (33-39) + trailing blank rows
# (40-44) edits described by the commit "2 Repos y arreglos a medias de
# pantallas. Revisar asociarEgryIng2".
#
# New shared-string cells are written in the exact order Excel would have
# assigned new shared-string table entries (first-use order), so the
# resulting xl/sharedStrings.xml unique-string ordering matches the
# target file exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "ENTREGA 5" section header (row 33) ------------------------------
# C33 becomes an "ENTREGA n" banner cell: copy the look of the existing
# "ENTREGA 4" banner (C23) so it picks up the same fill/border/bold style.
$ws.Range('C33').Value = 'ENTREGA 5'
$ws.Range('C23').Copy()
$ws.Range('C33').PasteSpecial(-4122)

# --- "ENTREGA 6" section header (row 36) -------------------------------
$ws.Range('C36').Value = 'ENTREGA 6'
$ws.Range('C23').Copy()
$ws.Range('C36').PasteSpecial(-4122)

# --- Open "DEFINIR PATRON DE INTERACION" TODO note (row 34) -----------
$ws.Range('C34').Value = 'DEFINIR PATRON DE INTERACION'
$ws.Range('C34').Font.Color = 255
$ws.Range('A34').Value = 'Otro'

# --- Row 37: Bitacora (Objetos) ---------------------------------------
$ws.Range('C37').Value = 'Hay una bitacora por usuario. Si este tiene proyectos asociados, podrá visualizarlos allí.'
$ws.Range('B37').Value = 'Bitacora'
$ws.Range('D37').Value = 'Que exista una unica y que permitamos a todos los usuarios acceder a toda la información.'
$ws.Range('E37').Value = 'Para mantener la seguridad de la información, más al ser un tema de dinero obtenido del estado.'
$ws.Range('A37').Value = 'Objetos'

# --- Row 38: Proyecto (Objetos) ----------------------------------------
$ws.Range('C38').Value = 'Existe un proyecto, el cual será registrado en una "OperacionRegitrada" para despues ser guardada en la bitacora.'
$ws.Range('C39').Value = 'La Bitacora es el unico elemento que se relaciona con la base de datos de Mongo db'
$ws.Range('B38').Value = 'Proyecto'
$ws.Range('D38').Value = 'Que toda la información este en una sola clase.'
$ws.Range('D39').Value = 'Que varias Clases interactuen con la base de datos.'
$ws.Range('E39').Value = 'Al ser la unica que interactua reducimos la posibilidad de fallas.'
$ws.Range('E38').Value = 'Menor acoplamiento.'
# C38's formatting was cleared to the default (no border) by the author.
$ws.Range('C38').Style = 'Normal'

$ws.Range('A38').Value = 'Objetos'

# --- Row 39: Bitacora (Persistencia) ------------------------------------
# Give the brand-new row 39 the same plain bordered style ("s=1") as the
# rest of the table before filling in its values.
$ws.Range('A35:E35').Copy()
$ws.Range('A39:E39').PasteSpecial(-4122)
$ws.Range('A39').Value = 'Persistencia'
$ws.Range('B39').Value = 'Bitacora'

# --- New trailing blank rows (40-44), matching the formatting of the
# existing blank row 35 (thin border, no fill) ---------------------------
$ws.Range('A40:E40').PasteSpecial(-4122)
$ws.Range('A41:E41').PasteSpecial(-4122)
$ws.Range('A42:E42').PasteSpecial(-4122)
$ws.Range('A43:E43').PasteSpecial(-4122)
$ws.Range('A44:E44').PasteSpecial(-4122)

# --- Column E got a bit wider to fit the new text, and the view scrolled
# down to the newly-entered rows ------------------------------------------
$ws.Columns('E').ColumnWidth = 173.28515625
$ws.Range('E38').Select()
$excel.ActiveWindow.Zoom = 85
